# This workbook's single data table (rows 2-32, columns A-R) is being
# reordered (re-sorted) - the weekly price rows are shuffled into a new
# row order while each row's own field values stay intact together.
#
# Approach: snapshot every data row (A:R) exactly as it is now, then
# write the rows back in the new order derived from the target edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 32
$lastCol  = 18   # column R

# Mapping: index i (0-based) corresponds to destination row ($firstRow + i);
# the value is the *original* (current) row number whose data should end up
# at that destination row once the reorder is complete.
$rowMap = @(17,24,6,31,7,11,22,23,14,25,4,16,30,15,9,5,13,12,26,27,28,32,18,8,19,29,2,3,10,20,21)

# 1) Snapshot all current rows (so overwriting destination rows during the
#    rewrite doesn't clobber data we still need to read later).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowValues = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowValues += ,$ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowValues
}

# 2) Write each destination row using the snapshotted source row's values.
for ($i = 0; $i -lt $rowMap.Length; $i++) {
    $destRow = $firstRow + $i
    $srcRow  = $rowMap[$i]
    $rowValues = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $rowValues[$c - 1]
    }
}
